# TN. RES values updated.
# Update the RES installed capacities (column C) from 25/25/30/30/30 to 50/50/50/50/50.
# All downstream formulas (Main!B7 total, and the Pg Winter/Summer S1-S3 profile
# sheets that VLOOKUP into 'RES installed') recalculate automatically.

$wb = $excel.ActiveWorkbook
$resSheet = $wb.Worksheets.Item("RES installed")

$resSheet.Range("C2").Value = 50
$resSheet.Range("C3").Value = 50
$resSheet.Range("C4").Value = 50
$resSheet.Range("C5").Value = 50
$resSheet.Range("C6").Value = 50

# Reflect the author's final UI state: "RES installed" sheet active/selected,
# with C7 as the selected cell.
$resSheet.Activate()
$resSheet.Range("C7").Select()
